$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (D2D - CPDMA)
$ws.Range("C3").Value = 0.30277
$ws.Range("D3").Value = 0.622094
$ws.Range("E3").Value = 1.273629
$ws.Range("F3").Value = 2.487014
$ws.Range("G3").Value = 4.664838
$ws.Range("H3").Value = 8.700206
$ws.Range("I3").Value = 14.958417
$ws.Range("J3").Value = 23.407006
$ws.Range("K3").Value = 31.685974
$ws.Range("L3").Value = 40.289131
$ws.Range("M3").Value = 42.780729
$ws.Range("N3").Value = 42.652847
$ws.Range("O3").Value = 44.723804
$ws.Range("P3").Value = 48.396205
$ws.Range("Q3").Value = 45.072236
$ws.Range("R3").Value = 45.223815
$ws.Range("S3").Value = 45.246745
$ws.Range("T3").Value = 45.144975

# Row 8 (H2D - CPDMA)
$ws.Range("C8").Value = 0.065431
$ws.Range("D8").Value = 0.132435
$ws.Range("E8").Value = 0.266465
$ws.Range("F8").Value = 0.516958
$ws.Range("G8").Value = 0.995807
$ws.Range("H8").Value = 1.858626
$ws.Range("I8").Value = 3.273398
$ws.Range("J8").Value = 4.763143
$ws.Range("K8").Value = 6.20474
$ws.Range("L8").Value = 8.658974000000001
$ws.Range("M8").Value = 10.332109
$ws.Range("N8").Value = 11.863943
$ws.Range("O8").Value = 12.335983
$ws.Range("P8").Value = 13.085151
$ws.Range("Q8").Value = 13.359447
$ws.Range("R8").Value = 13.33552
$ws.Range("S8").Value = 13.351706
$ws.Range("T8").Value = 13.457736

# Row 13 (D2H - CPDMA)
$ws.Range("C13").Value = 0.06561400000000001
$ws.Range("D13").Value = 0.132676
$ws.Range("E13").Value = 0.263979
$ws.Range("F13").Value = 0.521127
$ws.Range("G13").Value = 1.005211
$ws.Range("H13").Value = 1.871335
$ws.Range("I13").Value = 3.290276
$ws.Range("J13").Value = 4.756716
$ws.Range("K13").Value = 6.307795
$ws.Range("L13").Value = 8.823107
$ws.Range("M13").Value = 10.582898
$ws.Range("N13").Value = 11.807048
$ws.Range("O13").Value = 12.809033
$ws.Range("P13").Value = 13.124515
$ws.Range("Q13").Value = 13.649282
$ws.Range("R13").Value = 13.692114
$ws.Range("S13").Value = 13.660944
$ws.Range("T13").Value = 13.316356
